# Avverkningsanmälningar - refresh "Förändrad" (last-changed) snapshot date
# and append the newly reported case (A 46522-2023).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 446
$oldChanged  = 45192
$newChanged  = 45202

# 1. Column C ("Förändrad") holds the date the source report was refreshed —
#    it is the same value across every existing row, and moves from
#    2023-09-23 (45192) to 2023-10-03 (45202).
for ($r = 2; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, 3).Value2 = $newChanged
}

# 2. Row 446 picks up an explicit row height once a new row is appended
#    right after it.
$ws.Rows.Item($lastDataRow).RowHeight = 15

# 3. Append the new case as row 447.
$newRow = $lastDataRow + 1

$ws.Cells.Item($newRow, 1).Value2 = "A 46522-2023"

$ws.Cells.Item($newRow, 2).Value2 = 45197
$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 3).Value2 = $newChanged
$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 4).Value2 = "UPPSALA LÄN"
$ws.Cells.Item($newRow, 5).Value2 = "TIERP"

$ws.Cells.Item($newRow, 7).Value2 = 16.1

for ($c = 8; $c -le 17; $c++) {
    $ws.Cells.Item($newRow, $c).Value2 = 0
}

# Artnamn column keeps the same wrap-text style as every other row, even
# though this case has no species listed yet.
$ws.Cells.Item($newRow, 18).WrapText = $true
